$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 5 new rows before the old row 65 (summary row), pushing it down to row 70
# and creating blank rows 64-69 (with formatting inherited from row 63 above).
$ws.Rows("64:68").Insert()

# ---------------------------------------------------------------------------
# Row 64: new data row ("Peer reviewed Exposé" / "Exposé reviewen")
# ---------------------------------------------------------------------------
$ws.Range("A64").Value = 16
$ws.Range("B64").Value = "Konzeptuelles Design"
$ws.Range("C64").Value = "[SEMINAR]"
$ws.Range("D64").Value = "Peer reviewed Exposé"
$ws.Range("E64").Value = "Exposé reviewen"
$ws.Range("F64").Value = 44334
$ws.Range("G64").Value = 44338
$ws.Range("H64").Formula = "=ROUNDUP(((SUM(K64-J64)*24*60/60)/0.25),0)*0.25"
$ws.Range("J64").Value = 0.375
$ws.Range("K64").Value = 0.41666666666666669

# ---------------------------------------------------------------------------
# Row 65: new data row ("Umfrage" / "Umfrage auswerten")
# ---------------------------------------------------------------------------
$ws.Range("A65").Value = 14
$ws.Range("B65").Value = "Konzeptuelles Design"
$ws.Range("C65").Value = "[TASK]"
$ws.Range("D65").Value = "Umfrage"
$ws.Range("E65").Value = "Umfrage auswerten"
$ws.Range("F65").Value = 44332
$ws.Range("G65").Value = 44338
$ws.Range("I65").Formula = "=ROUNDUP(((SUM(K65-J65)*24*60/60)/0.25),0)*0.25"
$ws.Range("J65").Value = 0.41666666666666669
$ws.Range("K65").Value = 0.5

# ---------------------------------------------------------------------------
# Apply number formats (copy formats only, so values stay untouched) to match
# the existing style palette already used by surrounding rows.
# ---------------------------------------------------------------------------
$ws.Range("F61:G61").Copy()
$ws.Range("F64:G69").PasteSpecial(-4122)

$ws.Range("H61").Copy()
$ws.Range("H64").PasteSpecial(-4122)
$ws.Range("H66:H69").PasteSpecial(-4122)

$ws.Range("I61").Copy()
$ws.Range("I64:I69").PasteSpecial(-4122)

$ws.Range("J61").Copy()
$ws.Range("J64:J65").PasteSpecial(-4122)

$ws.Range("K61").Copy()
$ws.Range("K64:K65").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 63's prefix cell was part of the validation range and must stay excluded
# once the range below it grows; clearing/resetting its validation reproduces
# the exact split the author ended up with (C41:C62 C64:C69).
# ---------------------------------------------------------------------------
$ws.Range("C63").Validation.Delete()

# ---------------------------------------------------------------------------
# Reflect the author's final cursor position / scroll state in the view.
# ---------------------------------------------------------------------------
$ws.Range("D63").Select()

Write-Host "done"
